$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Date paragraph: "17." -> "19."
# ---------------------------------------------------------------------------
$datePara = $d.Paragraphs.Item(3)
$dateRange = $d.Range($datePara.Range.Start, $datePara.Range.Start + 3)
$dateRange.Find.Execute("17.", $true, $false, $false, $false, $false, $true, 1, $false, "19.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Table of random numbers: replace the 5x5 block of values.
#    NOTE: deliberately avoid $d.Tables / Table.Cell() here -- merely
#    materialising a Table/Cell reference leaves this host's Paragraphs
#    collection in a broken state (Paragraph.Range.Text starts coming back
#    empty for every paragraph). Plain whole-document Find/Replace is used
#    instead, which is safe because every one of these numbers is unique in
#    the document.
# ---------------------------------------------------------------------------
$oldValues = @(
    "-2.033", "0.617", "-0.363", "-0.008", "-0.082",
    "2.001", "-0.310", "1.921", "0.621", "-1.432",
    "-0.675", "1.274", "0.380", "-2.311", "-0.300",
    "1.424", "-2.047", "0.502", "0.838", "-1.618",
    "-0.224", "1.627", "-2.426", "0.409", "0.060"
)

$newValues = @(
    "-0.829", "-0.573", "1.527", "-0.033", "-1.338",
    "-0.134", "-0.095", "0.602", "-1.368", "0.078",
    "0.145", "1.031", "0.448", "-1.350", "-0.105",
    "1.191", "-1.234", "-0.832", "0.163", "-1.205",
    "1.058", "1.366", "0.694", "0.188", "1.525"
)

for ($i = 0; $i -lt $oldValues.Count; $i++) {
    $d.Content.Find.Execute($oldValues[$i], $true, $false, $false, $false, $false, $true, 1, $false, $newValues[$i], 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Insert new Heading 6 paragraph before "And here is a bunch of
#    informative text.", and switch that paragraph's style to
#    "First Paragraph".
# ---------------------------------------------------------------------------
$paras = @($d.Paragraphs)
$infoParaIndex = -1
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text -like "*And here is a bunch of informative text.*") {
        $infoParaIndex = $i
        break
    }
}

$infoPara = $paras[$infoParaIndex]
$heading6Text = "l I have defined the 6th level header with a small font and white ink"
$insertPoint = $d.Range($infoPara.Range.Start, $infoPara.Range.Start)
$insertPoint.InsertBefore($heading6Text + "`r")

$paras = @($d.Paragraphs)
$heading6Para = $paras[$infoParaIndex]
$heading6Para.Range.Style = "Heading 6"
$heading6Range = $d.Range($heading6Para.Range.Start, $heading6Para.Range.Start + $heading6Text.Length)
$d.Bookmarks.Add("Xd8d52b0cce2a3ef2e3573ad4e5a5854bc2429fd", $heading6Range)

$paras = @($d.Paragraphs)
$infoPara = $paras[$infoParaIndex + 1]
$infoPara.Range.Style = "First Paragraph"

# ---------------------------------------------------------------------------
# 4. Insert new Heading 5 paragraph before "Yet nore text introducing the
#    following equation.", and switch that paragraph's style to
#    "First Paragraph".
# ---------------------------------------------------------------------------
$paras = @($d.Paragraphs)
$yetParaIndex = -1
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text -like "*Yet nore text introducing the following equation.*") {
        $yetParaIndex = $i
        break
    }
}

$yetPara = $paras[$yetParaIndex]
$heading5Text = "p and the 5th level header as requiring a page break beforehand"
$insertPoint = $d.Range($yetPara.Range.Start, $yetPara.Range.Start)
$insertPoint.InsertBefore($heading5Text + "`r")

$paras = @($d.Paragraphs)
$heading5Para = $paras[$yetParaIndex]
$heading5Para.Range.Style = "Heading 5"
$heading5Range = $d.Range($heading5Para.Range.Start, $heading5Para.Range.Start + $heading5Text.Length)
$d.Bookmarks.Add("X61634945cc9478cc10a632f4762567517962fb3", $heading5Range)

$paras = @($d.Paragraphs)
$yetPara = $paras[$yetParaIndex + 1]
$yetPara.Range.Style = "First Paragraph"

# ---------------------------------------------------------------------------
# 5. Append a new BodyText paragraph at the very end of the document.
# ---------------------------------------------------------------------------
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertBefore("`rOf course you will need to provide your own reference document and modify the styles within that to get the formatting you desire.")

$paras = @($d.Paragraphs)
$lastPara = $paras[$paras.Count - 1]
$lastPara.Range.Style = "Body Text"
